$wb = $excel.ActiveWorkbook

# --- Update the "Last Updated" timestamp on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 09:33 AM"

# --- Insert a new row 2 on the "Stock List" sheet with the new data ---
$stock = $wb.Worksheets.Item("Stock List")
$stock.Rows.Item(2).Insert()
$stock.Rows.Item(2).ClearFormats()

$stock.Cells.Item(2, 1).Value = "📋"
$stock.Cells.Item(2, 2).Value = "CAPTRU-RE1"
$stock.Cells.Item(2, 3).Value = "CAPTRU-RE1"
$stock.Cells.Item(2, 4).Value = 5.67
$stock.Cells.Item(2, 5).Value = -11.9565
$stock.Cells.Item(2, 6).Value = "N/A"
$stock.Cells.Item(2, 7).Value = "N/A"
$stock.Cells.Item(2, 8).Value = 0

# Remove the now-duplicated last row (the insert pushed it past the original range)
$stock.Rows.Item(77).Delete()
